$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap Starting / Ending SoC (%) values (row 6 / row 7) ---
$ws.Range("B6").Value = 99
$ws.Range("B7").Value = 14

# --- Relabel rows with units, fix wording ---
$ws.Range("A8").Value  = "Total distance covered (km)"
$ws.Range("A9").Value  = "Total energy consumption(WH/KM)"
$ws.Range("A10").Value = "Total SOC consumed(%)"

$ws.Range("A12").Value = "Peak Power(kW)"
$ws.Range("A13").Value = "Average Power(kW)"
$ws.Range("A14").Value = "Total Energy Regenerated(kWh)"

# --- Regenerative Effectiveness: relabel and flip sign (anomaly fix) ---
$ws.Range("A15").Value = "Regenerative Effectiveness(%)"
$ws.Range("B15").Value = 5.228834716066615

# --- Swap Lowest/Highest Cell Voltage rows (row 16 / row 17) ---
$ws.Range("A16").Value = "Highest Cell Voltage(V)"
$ws.Range("B16").Value = 3.34
$ws.Range("A17").Value = "Lowest Cell Voltage(V)"
$ws.Range("B17").Value = 3.107

$ws.Range("A18").Value = "Difference in Cell Voltage(V)"
$ws.Range("A19").Value = "Minimum Temperature(C)"
$ws.Range("A20").Value = "Maximum Temperature(C)"

# --- Difference in Temperature: relabel and fill in missing value ---
$ws.Range("A21").Value = "Difference in Temperature(C)"
$ws.Range("B21").Value = 7

$ws.Range("A22").Value = "Maximum Fet Temperature-BMS(C)"
$ws.Range("A23").Value = "Maximum Afe Temperature-BMS(C)"
$ws.Range("A24").Value = "Maximum PCB Temperature-BMS(C)"
$ws.Range("A25").Value = "Maximum MCU Temperature(C)"
$ws.Range("A26").Value = "Maximum Motor Temperature(C)"
$ws.Range("A27").Value = "Abnormal Motor Temperature Detected(C)"

# --- Swap lowest/highest cell temp rows (row 28 / row 29) ---
$ws.Range("A28").Value = "highest cell temp(C)"
$ws.Range("B28").Value = 42
$ws.Range("A29").Value = "lowest cell temp(C)"
$ws.Range("B29").Value = 31

$ws.Range("A30").Value = "Difference between Highest and Lowest Cell Temperature at 100% SOC(C)"

# --- Remove "Maximum BMS Temperature in C" row: everything below shifts up one row,
#     and a brand-new "Time spent in 80-90 km/h" row is appended at the end ---
$ws.Range("A31").Value = "Battery Voltage(V)"
$ws.Range("B31").Value = 55

$ws.Range("A32").Value = "Total energy charged(kWh)"
$ws.Range("B32").Value = 1.851597901388889

$ws.Range("A33").Value = "Electricity consumption units(kW)"
$ws.Range("B33").Value = 0.00000006182627123281674

$ws.Range("A34").Value = "Idling time percentage"
$ws.Range("B34").Value = 9.657230379211999

$ws.Range("A35").Value = "Time spent in 0-10 km/h"
$ws.Range("B35").Value = 9.07494454420619

$ws.Range("A36").Value = "Time spent in 10-20 km/h"
$ws.Range("B36").Value = 10.38475757895849

$ws.Range("A37").Value = "Time spent in 20-30 km/h"
$ws.Range("B37").Value = 21.83505862469631

$ws.Range("A38").Value = "Time spent in 30-40 km/h"
$ws.Range("B38").Value = 41.82291116509982

$ws.Range("A39").Value = "Time spent in 40-50 km/h"
$ws.Range("B39").Value = 6.743160452096757

$ws.Range("A40").Value = "Time spent in 50-60 km/h"
$ws.Range("B40").Value = 0.05809654589627126

$ws.Range("A41").Value = "Time spent in 60-70 km/h"
$ws.Range("B41").Value = 0.0250871448188444

$ws.Range("A42").Value = "Time spent in 70-80 km/h"
$ws.Range("B42").Value = 0

$ws.Range("A43").Value = "Time spent in 80-90 km/h"
$ws.Range("B43").Value = 0
